$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 44.50020533333333
$ws.Range("N2").Value = 133.500616
$ws.Range("O2").Value = 0.2926972930209797
$ws.Range("P2").Value = 0.2926972930209797
$ws.Range("Q2").Value = 4286.292974859845
$ws.Range("R2").Value = 38576.63677373861
$ws.Range("S2").Value = 0.1115125350497709
$ws.Range("T2").Value = 0.1115125350497709
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 28.185334
$ws.Range("N3").Value = 84.55600199999999
$ws.Range("O3").Value = 0.1853872561462678
$ws.Range("P3").Value = 0.1853872561462678
$ws.Range("Q3").Value = 2714.832397139164
$ws.Range("R3").Value = 24433.49157425248
$ws.Range("S3").Value = 0.07062929310149027
$ws.Range("T3").Value = 0.07062929310149026
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 62.31760433333333
$ws.Range("N4").Value = 186.952813
$ws.Range("O4").Value = 0.4098901108273345
$ws.Range("P4").Value = 0.4098901108273344
$ws.Range("Q4").Value = 6002.478138319499
$ws.Range("R4").Value = 54022.3032448755
$ws.Range("S4").Value = 0.1561609431998109
$ws.Range("T4").Value = 0.1561609431998109
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.031762
$ws.Range("N5").Value = 51.09528599999999
$ws.Range("O5").Value = 0.1120253400054181
$ws.Range("P5").Value = 0.1120253400054181
$ws.Range("Q5").Value = 1640.512021534452
$ws.Range("R5").Value = 14764.60819381007
$ws.Range("S5").Value = 0.04267968973980667
$ws.Range("T5").Value = 0.04267968973980666
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.50020533333333
$ws.Range("N6").Value = 133.500616
$ws.Range("O6").Value = 0.2926972930209797
$ws.Range("P6").Value = 0.2926972930209797
$ws.Range("Q6").Value = 821.6818437466684
$ws.Range("R6").Value = 7395.136593720016
$ws.Range("S6").Value = 0.02137693945280464
$ws.Range("T6").Value = 0.02137693945280464
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.185334
$ws.Range("N7").Value = 84.55600199999999
$ws.Range("O7").Value = 0.1853872561462678
$ws.Range("P7").Value = 0.1853872561462678
$ws.Range("Q7").Value = 520.4330414715613
$ws.Range("R7").Value = 4683.897373244051
$ws.Range("S7").Value = 0.01353962692670442
$ws.Range("T7").Value = 0.01353962692670442
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 62.31760433333333
$ws.Range("N8").Value = 186.952813
$ws.Range("O8").Value = 0.4098901108273345
$ws.Range("P8").Value = 0.4098901108273344
$ws.Range("Q8").Value = 1150.674331566126
$ws.Range("R8").Value = 10356.06898409514
$ws.Range("S8").Value = 0.02993603388341299
$ws.Range("T8").Value = 0.02993603388341299
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.031762
$ws.Range("N9").Value = 51.09528599999999
$ws.Range("O9").Value = 0.1120253400054181
$ws.Range("P9").Value = 0.1120253400054181
$ws.Range("Q9").Value = 314.4859556846039
$ws.Range("R9").Value = 2830.373601161436
$ws.Range("S9").Value = 0.008181691349991495
$ws.Range("T9").Value = 0.008181691349991494
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 44.50020533333333
$ws.Range("N10").Value = 133.500616
$ws.Range("O10").Value = 0.2926972930209797
$ws.Range("P10").Value = 0.2926972930209797
$ws.Range("Q10").Value = 5469.493136894951
$ws.Range("R10").Value = 49225.43823205457
$ws.Range("S10").Value = 0.1422947634960541
$ws.Range("T10").Value = 0.1422947634960541
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 28.185334
$ws.Range("N11").Value = 84.55600199999999
$ws.Range("O11").Value = 0.1853872561462678
$ws.Range("P11").Value = 0.1853872561462678
$ws.Range("Q11").Value = 3464.242237071593
$ws.Range("R11").Value = 31178.18013364434
$ws.Range("S11").Value = 0.09012599842057566
$ws.Range("T11").Value = 0.09012599842057564
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 62.31760433333333
$ws.Range("N12").Value = 186.952813
$ws.Range("O12").Value = 0.4098901108273345
$ws.Range("P12").Value = 0.4098901108273344
$ws.Range("Q12").Value = 7659.41879718896
$ws.Range("R12").Value = 68934.76917470065
$ws.Range("S12").Value = 0.1992680416602499
$ws.Range("T12").Value = 0.1992680416602499
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.031762
$ws.Range("N13").Value = 51.09528599999999
$ws.Range("O13").Value = 0.1120253400054181
$ws.Range("P13").Value = 0.1120253400054181
$ws.Range("Q13").Value = 2093.363495076941
$ws.Range("R13").Value = 18840.27145569247
$ws.Range("S13").Value = 0.05446110928157248
$ws.Range("T13").Value = 0.05446110928157247
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 44.50020533333333
$ws.Range("N14").Value = 133.500616
$ws.Range("O14").Value = 0.2926972930209797
$ws.Range("P14").Value = 0.2926972930209797
$ws.Range("Q14").Value = 673.1627496149093
$ws.Range("R14").Value = 6058.464746534184
$ws.Range("S14").Value = 0.01751305502235
$ws.Range("T14").Value = 0.01751305502235
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 28.185334
$ws.Range("N15").Value = 84.55600199999999
$ws.Range("O15").Value = 0.1853872561462678
$ws.Range("P15").Value = 0.1853872561462678
$ws.Range("Q15").Value = 426.364705334122
$ws.Range("R15").Value = 3837.282348007097
$ws.Range("S15").Value = 0.01109233769749749
$ws.Range("T15").Value = 0.01109233769749749
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 62.31760433333333
$ws.Range("N16").Value = 186.952813
$ws.Range("O16").Value = 0.4098901108273345
$ws.Range("P16").Value = 0.4098901108273344
$ws.Range("Q16").Value = 942.6898048719263
$ws.Range("R16").Value = 8484.208243847337
$ws.Range("S16").Value = 0.02452509208386058
$ws.Range("T16").Value = 0.02452509208386058
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.031762
$ws.Range("N17").Value = 51.09528599999999
$ws.Range("O17").Value = 0.1120253400054181
$ws.Range("P17").Value = 0.1120253400054181
$ws.Range("Q17").Value = 257.6425805864459
$ws.Range("R17").Value = 2318.783225278014
$ws.Range("S17").Value = 0.006702849634047455
$ws.Range("T17").Value = 0.006702849634047454
